# Apply the "new standards" text edits to both ScenarioA and ScenarioB sheets.
# - Row 5 (description row), column E: simplify "Which package this node belongs to"
#   to "Which package this belongs to".
# - Row 6 (example row), columns D, E, F: change placeholder "-" to
#   "Scenario-dependent" for pWeight_k's dataPackage/dataSource columns.
# - Row 5's height shrinks from 45 to 30 because the new text is shorter.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("E5").Value2 = "Which package this belongs to"

    $ws.Range("D6").Value2 = "Scenario-dependent"
    $ws.Range("E6").Value2 = "Scenario-dependent"
    $ws.Range("F6").Value2 = "Scenario-dependent"

    $ws.Rows.Item(5).RowHeight = 30
}
